$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("A26").Value = 'x_4*x_5'
$ws.Range("B26").Value = 'Did not compute.'
$ws.Range("C26").Value = 'Did not compute.'
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 'Did not compute.'
$ws.Range("F26").Value = 'F_std is non-homogeneous or non-cubic; abort.'
$ws.Range("G26").Value = 'Did not compute.'

# Row 27
$ws.Range("A27").Value = 'x_3*x_4*x_5'
$ws.Range("B27").Value = 'Did not compute.'
$ws.Range("C27").Value = 'Did not compute.'
$ws.Range("D27").Value = 11
$ws.Range("E27").Value = 'Did not compute.'
$ws.Range("F27").Value = 'X contains a possibly relevant geometric plane; abort.'
$ws.Range("G27").Value = 'Did not compute.'

# Row 28
$ws.Range("A28").Value = 'x_1^3+x_2^3+x_3*x_4*x_5'
$ws.Range("B28").Value = 'Did not compute.'
$ws.Range("C28").Value = 'Did not compute.'
$ws.Range("D28").Value = 11
$ws.Range("E28").Value = 'Did not compute.'
$ws.Range("F28").Value = 'X contains a possibly relevant geometric plane; abort.'
$ws.Range("G28").Value = 'Did not compute.'

# Row 29
$ws.Range("A29").Value = 'x_1^3+x_2^3+x_3^3+x_3*x_4*x_5'
$ws.Range("B29").Value = 'Did not compute.'
$ws.Range("C29").Value = 'Did not compute.'
$ws.Range("D29").Value = 11
$ws.Range("E29").Value = 'Did not compute.'
$ws.Range("F29").Value = 'X contains a possibly relevant geometric plane; abort.'
$ws.Range("G29").Value = 'Did not compute.'

# Row 30
$ws.Range("A30").Value = '(2*x_1^3 + x_1*x_2*x_3 + 2*x_2^3 + x_3^3) + 2*(x_1^2+x_2*x_3)*x_4 + 2*(x_2^2+x_1*x_3)*x_5 + 2*x_3*x_4*x_5'
$ws.Range("B30").Value = 'Did not compute.'
$ws.Range("C30").Value = 'Did not compute.'
$ws.Range("D30").Value = 11
$ws.Range("E30").Value = 110
$ws.Range("F30").Value = 110
$ws.Range("G30").Value = '[(0 : 0 : 0 : 1 : 0), (0 : 0 : 0 : 0 : 1)]'

# Row 31
$ws.Range("A31").Value = '-x_1^2*x_2 - 2*x_2^3 - 2*x_1^2*x_3 + x_1*x_2*x_3 - x_2^2*x_3 - x_1*x_3^2 - 2*x_2*x_3^2 + 2*x_2^2*x_4 - x_1*x_3*x_4 - 2*x_2*x_3*x_4 - 2*x_3^2*x_4 + x_1*x_4^2 + x_3*x_4^2 + x_1^2*x_5 + x_1*x_2*x_5 + 2*x_1*x_3*x_5 + x_3^2*x_5 - 2*x_1*x_4*x_5 + x_2*x_4*x_5 - 2*x_3*x_4*x_5 + x_1*x_5^2 - x_2*x_5^2 + x_3*x_5^2'
$ws.Range("B31").Value = 'Did not compute.'
$ws.Range("C31").Value = 'Did not compute.'
$ws.Range("D31").Value = 101
$ws.Range("E31").Value = 'Did not compute.'
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = '[(0 : 0 : 0 : 1 : 1)]'

# Row 32
$ws.Range("A32").Value = '-x_1^2*x_2 - 2*x_2^3 - 2*x_1^2*x_3 + x_1*x_2*x_3 - x_2^2*x_3 - x_1*x_3^2 - 2*x_2*x_3^2 + 2*x_2^2*x_4 - x_1*x_3*x_4 - 2*x_2*x_3*x_4 - 2*x_3^2*x_4 + x_1*x_4^2 + x_3*x_4^2 + x_1^2*x_5 + x_1*x_2*x_5 + 2*x_1*x_3*x_5 + x_3^2*x_5 - 2*x_1*x_4*x_5 + x_2*x_4*x_5 - 2*x_3*x_4*x_5 + x_1*x_5^2 - x_2*x_5^2 + x_3*x_5^3'
$ws.Range("B32").Value = 'Did not compute.'
$ws.Range("C32").Value = 'Did not compute.'
$ws.Range("D32").Value = 10007
$ws.Range("E32").Value = 'Did not compute.'
$ws.Range("F32").Value = -20014
$ws.Range("G32").Value = '[(0 : 0 : 0 : 1 : 1)]'

# Row 33
$ws.Range("A33").Value = '-x_1^2*x_2 - 2*x_2^3 - 2*x_1^2*x_3 + x_1*x_2*x_3 - x_2^2*x_3 - x_1*x_3^2 - 2*x_2*x_3^2 + 2*x_2^2*x_4 - x_1*x_3*x_4 - 2*x_2*x_3*x_4 - 2*x_3^2*x_4 + x_1*x_4^2 + x_3*x_4^2 + x_1^2*x_5 + x_1*x_2*x_5 + 2*x_1*x_3*x_5 + x_3^2*x_5 - 2*x_1*x_4*x_5 + x_2*x_4*x_5 - 2*x_3*x_4*x_5 + x_1*x_5^2 - x_2*x_5^2 + x_3*x_5^4'
$ws.Range("B33").Value = 'Did not compute.'
$ws.Range("C33").Value = 'Did not compute.'
$ws.Range("D33").Value = 10009
$ws.Range("E33").Value = 'Did not compute.'
$ws.Range("F33").Value = 380342
$ws.Range("G33").Value = '[(0 : 0 : 0 : 1 : 1)]'

# Row 34
$ws.Range("A34").Value = '(2*x_1^3 + x_1*x_2*x_3 + 2*x_2^3 + x_3^3) + 2*(x_1^2+x_2*x_3)*x_4 + 2*(x_2^2+x_1*x_3)*x_5 + 2*x_3*x_4*x_5'
$ws.Range("B34").Value = 'Did not compute.'
$ws.Range("C34").Value = 'Did not compute.'
$ws.Range("D34").Value = '[101,10007,10009]'
$ws.Range("E34").Value = 'Did not compute.'
$ws.Range("F34").Value = '[10100, 100130042, 100170072]'
$ws.Range("G34").Value = '[(0 : 0 : 0 : 1 : 0), (0 : 0 : 0 : 0 : 1)]'

# A31:A33 use the "text" number format (same style as A24/A25 above them)
$ws.Range("A31:A33").NumberFormat = "@"

# Column D widens to fit the new longer numeric entries (11, 101, 10007, 10009, and the [101,10007,10009] text)
$ws.Columns("D").ColumnWidth = 11.953125

# Update the view to match where the user left off after adding the new rows
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("G39").Select()
